$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.175204038619995
$ws.Range("B1").Value = 2.200358867645264
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.368507623672485
$ws.Range("E1").Value = 1.227415084838867
